$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: "MSAL is a primarily mobile website ... to other. MSAL has 4 main
# pages:" -- merge the "to" run (and its surrounding gramStart/gramEnd
# proofErr markers) back into a single run.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found = $r1.Find.Execute("MSAL is a primarily mobile website where you can list your all-time favourite Spotify artist and show it to other. MSAL has 4 main pages:")
if (-not $found) { throw "Change1: paragraph not found" }
$p1 = $r1.Paragraphs(1).Range
$p1.InsertXML("<w:p $w><w:r><w:tab/><w:t>MSAL is a primarily mobile website where you can list your all-time favourite Spotify artist and show it to other. MSAL has 4 main pages:</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Change 2: "Search Page: ... a spotify artist and you can add it ..." --
# merge the "artist" run (gramStart/gramEnd) into the following run.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found = $r2.Find.Execute("Search Page: Where you can search the name of a spotify artist and you can add it to your favourite page or see the artist details")
if (-not $found) { throw "Change2: paragraph not found" }
$p2 = $r2.Paragraphs(1).Range
$p2.InsertXML("<w:p $w><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t xml:space=""preserve"">Search Page: Where you can search the name of a </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>spotify</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> artist and you can add it to your favourite page or see the artist details</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Change 3: "Detail Page: See more information about the artist name, ..." --
# merge the "artist" run (gramStart/gramEnd) and the trailing space run.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found = $r3.Find.Execute("Detail Page: See more information about the artist name, and discography (album and tracks)")
if (-not $found) { throw "Change3: paragraph not found" }
$p3 = $r3.Paragraphs(1).Range
$p3.InsertXML("<w:p $w><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t xml:space=""preserve"">Detail Page: See </w:t></w:r><w:r><w:t>more</w:t></w:r><w:r><w:t xml:space=""preserve""> information about the artist </w:t></w:r><w:r><w:t>name, and discography (album and tracks)</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Change 4: "Favorite page is the page that is used to list all ... Favorite
# page get and store the data from the local storage, ..." -- merge the
# "page" run (gramStart/gramEnd).
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found = $r4.Find.Execute("is the page that is used to list all")
if (-not $found) { throw "Change4: paragraph not found" }
$p4 = $r4.Paragraphs(1).Range
$p4.InsertXML("<w:p $w><w:r><w:tab/></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>Favorite</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> page is the page that is used to list all </w:t></w:r><w:r><w:t>user’s</w:t></w:r><w:r><w:t xml:space=""preserve""> </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>favorite</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> artist (marked by the yellow star). </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>Favorite</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> page </w:t></w:r><w:r><w:t xml:space=""preserve"">get and store the data from the local storage, for </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>everytime</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> the user click the star (</w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>favorite</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve"">) button. You can also remove unwanted </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>favorite</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> artist by clicking the star button again, it should disappear from the list.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Change 5: "As per usual, user can add it to their favourite, or see the
# artist detail" -- merge the "detail" run (gramStart/gramEnd).
# ---------------------------------------------------------------------------
$r5 = $d.Content
$found = $r5.Find.Execute("As per usual, user can add it to their favourite, or see the artist detail")
if (-not $found) { throw "Change5: paragraph not found" }
$p5 = $r5.Paragraphs(1).Range
$p5.InsertXML("<w:p $w><w:pPr><w:ind w:firstLine=""720""/></w:pPr><w:r><w:t>As per usual, user can add it to their favourite, or see the artist detail</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Change 6: "Detail page have the URL ..." -- merge the "page" run
# (gramStart/gramEnd).
# ---------------------------------------------------------------------------
$r6 = $d.Content
$found = $r6.Find.Execute("Detail page have the URL")
if (-not $found) { throw "Change6: paragraph not found" }
$p6 = $r6.Paragraphs(1).Range
$p6.InsertXML("<w:p $w><w:r><w:tab/><w:t>Detail page have the URL " + [char]8220 + "/detail/{artist-name}" + [char]8221 + ", so user need to click the specified artist to get the right detail page. It consists of artists info, from their image, name, and discography (albums and tracks inside it). The favourite button (star) is also provided for easy access.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Change 7: drop the trailing "BONUS: Light and Dark mode" section (page
# break paragraph + heading paragraph + body paragraph) and replace it with
# a single empty paragraph.
# ---------------------------------------------------------------------------
$r7a = $d.Content
$found = $r7a.Find.Execute("Using your browser setting, you can also view the page in light mode")
if (-not $found) { throw "Change7: paragraph not found" }
$endPara = $r7a.Paragraphs(1)

$r7b = $d.Content
$found2 = $r7b.Find.Execute("BONUS: Light and Dark mode")
if (-not $found2) { throw "Change7: heading not found" }
$bonusPara = $r7b.Paragraphs(1)
$breakPara = $bonusPara.Previous()

$breakRange = $d.Range($breakPara.Range.Start, $endPara.Range.End)
$breakRange.InsertXML("<w:p $w/>")
